$d = $word.ActiveDocument

# Change 1: Question 10 answer - "Q-Q Plot" -> "histogram"
$rng1 = $d.Content
$rng1.Find.Execute(
    "Q-Q Plot will help us know if the data are normally distributed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "histogram will help us know if the data are normally distributed.",
    2)

# Change 2: Part (b) explanation about non-certified / certified cars
$rng2 = $d.Content
$oldText2 = "The sample size for non-certified cars is large, so we can assume that it is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution, but a Q-Q Plot will help us know more confidently."
$newText2 = "The sample size for non-certified cars is large, so we can assume that the sampling distribution of the sample mean is normally distributed. With a sample size of n = 24 for certified cars, it is likely large enough to assume a normal distribution for the mean, but a histogram of the data will help us know if the data is normally distributed. If it is, then the distribution of the sample mean will also be normally distributed."
$rng2.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $newText2, 2)
